# Wrap the "CVPR 23/24, ICCV 23, ICML 22, ICRA 21/22" reviewer list with a
# pair of co-located OLE_LINK bookmarks (OLE_LINK1 / OLE_LINK2), as produced
# by copying that text from an external OLE source (e.g. Excel/another Word
# doc) and pasting it into this CV.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("CVPR 23/24, ICCV 23, ICML 22, ICRA 21/22", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $d.Bookmarks.Add("OLE_LINK1", $rng)
    $d.Bookmarks.Add("OLE_LINK2", $rng)
}
